# Add a new worksheet "11_" to the workbook, as a duplicate of the last
# existing sheet ("10_"), then overwrite its question/answer content with
# the new "insulin in the tissue" (dX/dt) question. Mirrors an
# Excel UI workflow of right-click > "Move or Copy..." > Create a copy,
# then editing the new sheet's cells.

$wb = $excel.ActiveWorkbook

# The template sheet ("10_") is the last sheet in the workbook.
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate it, placing the copy immediately after the template sheet.
$templateSheet.Copy([System.Reflection.Missing]::Value, $templateSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "11_"

# --- Update the new sheet's content -----------------------------------
# Row 1: the prompt/question (column A), "Answer" header (B1), "Definitions"
# header (C1) are inherited from the template; only A1's text changes.
$newSheet.Range("A1").Value = "'X is the level of insulin in the tissue at any time.  Look at the dX/dt equation, and as in the previous question, match each term with the possible effects of the term. "

# Row 2: term (A2), correct-letter flag (B2, unchanged = C), description
# (C2, unchanged), explanation (D2) change.
$newSheet.Range("A2").Value = "'k3 (I(t) - I_b)"
$newSheet.Range("D2").Value = "Yes!  'I' can either be larger or smaller than I_b, so this term can be negative or positive."

# Row 3: term (A3), correct-letter flag (B3, unchanged = B), description
# (C3, unchanged), explanation (D3) change.
$newSheet.Range("A3").Value = "'- k2 X(t)"
$newSheet.Range("D3").Value = "Yes!  Since  X is always positive, this term will always be negative."

# Row heights: rows 1 and 3 differ from the template's; row 2/4/5 match.
$newSheet.Rows.Item(1).RowHeight = 60
$newSheet.Rows.Item(3).RowHeight = 45

# Selection / scroll state on the new sheet.
$newSheet.Activate()
$newSheet.Range("C10:C11").Select()
